$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 16.609
$ws.Range("D10").Value = -8.402999999999999
$ws.Range("D12").Value = -7.290000000000001
$ws.Range("E15").Value = 16.361
$ws.Range("D18").Value = -8.318999999999999
$ws.Range("E20").Value = 16.396
$ws.Range("E29").Value = 17.05
$ws.Range("E30").Value = 16.578
$ws.Range("E31").Value = 16.568
$ws.Range("D37").Value = -8.081999999999999
$ws.Range("E40").Value = 16.627
$ws.Range("D55").Value = -8.318999999999999
$ws.Range("D68").Value = -7.188
$ws.Range("E68").Value = 17.777
$ws.Range("E76").Value = 16.558
$ws.Range("D77").Value = -7.840000000000001
$ws.Range("D78").Value = -8.08
$ws.Range("E87").Value = 16.436
$ws.Range("E88").Value = 16.332
$ws.Range("E96").Value = 16.325
$ws.Range("E98").Value = 16.299
$ws.Range("E101").Value = 16.625
$ws.Range("E102").Value = 16.617
